$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (col F) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 264
$wsExhibit.Range("F4").Value = 910
$wsExhibit.Range("F6").Value = 46

# Sheet "全部类型" (sheet4): update 想去人数 (col F) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 264
$wsAll.Range("F5").Value = 910
$wsAll.Range("F7").Value = 46
